$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns remain stored as text (not auto-converted numbers/dates)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.868.43"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "2.292.50"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "300.69"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "98.97"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.510"
$ws.Range("E9").Value = "  +3.32%  "
$ws.Range("D10").Value = "35.95"
$ws.Range("E10").Value = "  +8.02%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.116"
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "18.28"
$ws.Range("E13").Value = "  +9.10%  "
$ws.Range("D14").Value = "6.92"
$ws.Range("E14").Value = "  +2.48%  "
$ws.Range("D15").Value = "2.651.49"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "2.235.35"
$ws.Range("E16").Value = "  -2.74%  "
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "42.778.61"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "12.45"
$ws.Range("E19").Value = "  +8.29%  "
$ws.Range("D20").Value = "0.0₃0901"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("E21").Value = "  +1.78%  "
$ws.Range("D22").Value = "67.64"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").Value = "235.09"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  +11.44%  "
$ws.Range("D25").Value = "1.01"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").Value = "24.93"
$ws.Range("E27").Value = "  +2.81%  "
$ws.Range("E28").Value = "  +14.88%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "34.35"
$ws.Range("E29").Value = "  +1.98%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "166.87"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "4.99"
$ws.Range("E33").Value = "  +1.72%  "
$ws.Range("E34").Value = "  +5.71%  "
$ws.Range("D35").Value = "4.65"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("E40").Value = "  +2.35%  "
$ws.Range("D41").Value = "0.109"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("D43").Value = "0.0291"
$ws.Range("E43").Value = "  +4.26%  "
$ws.Range("D44").Value = "1.971.87"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").Value = "10.12"
$ws.Range("E45").Value = "  +3.58%  "
$ws.Range("D46").Value = "17.51"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("D48").Value = "55.15"
$ws.Range("E48").Value = "  +4.80%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "1.54"
$ws.Range("E49").Value = "  +3.97%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.519.13"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").Value = "70.58"
$ws.Range("E51").Value = "  +1.43%  "
